# Edit: "ELSE IF b > a AND b > c THEN" -> "ELSE IF b > c THEN"
# The author selected "a AND b > " (between "ELSE IF b > " and "c THEN") and
# deleted it. This leaves the insertion point (and Word's "_GoBack" last-edit
# bookmark) sitting right before "THEN", splitting the paragraph's single run
# into three runs around that point, and the previously-existing "_GoBack"
# bookmark (at the end of the document, after the inline picture) is removed
# since a document can only have one bookmark of a given name.

$d = $word.ActiveDocument

$before = "ELSE IF b > "
$middleRemoved = "a AND b > "
$afterKeep = "c "
$tail = "THEN"

$full = $d.Content.Text
$idx = $full.IndexOf($before + $middleRemoved + $afterKeep + $tail)
if ($idx -lt 0) {
    throw "Could not locate target paragraph text"
}

# 1) Delete "a AND b > "
$delStart = $idx + $before.Length
$delEnd = $delStart + $middleRemoved.Length
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete()

# Re-resolve the (now shorter) text position.
$full2 = $d.Content.Text
$idx2 = $full2.IndexOf($before + $afterKeep + $tail)
if ($idx2 -lt 0) {
    throw "Could not locate edited paragraph text"
}

# 2) Force a run split after "ELSE IF b > " by touching (no-op) formatting.
$r1End = $idx2 + $before.Length
$r1 = $d.Range($idx2, $r1End)
$r1.Font.Bold = 1
$r1.Font.Bold = 0

# 3) Force a run split after "c " (i.e. before "THEN").
$r2End = $r1End + $afterKeep.Length
$r2 = $d.Range($r1End, $r2End)
$r2.Font.Bold = 1
$r2.Font.Bold = 0

# 4) Move the "_GoBack" bookmark to sit right between "c " and "THEN" --
#    Bookmarks.Add with an existing name relocates it (removing the old one).
$bmRange = $d.Range($r2End, $r2End)
$d.Bookmarks.Add("_GoBack", $bmRange)
